$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEAVE CREDITS")
$cv = $wb.Worksheets.Item("CONVERTION")

# --- Header: set the "period covered" start date (F3, merged F3:G3) ---
$ws.Range("F3").Value = 45110

# --- VL/SL earned entries for Jul-Nov 2023 (rows 11-15) ---
$ws.Range("A11").Value = 45110
$ws.Range("C11").Value = 1.167

$ws.Range("A12").Value = 45169
$ws.Range("C12").Value = 1.25

$ws.Range("A13").Value = 45199
$ws.Range("C13").Value = 1.25

$ws.Range("A14").Value = 45230
$ws.Range("C14").Value = 1.25

$ws.Range("A15").Value = 45260
$ws.Range("C15").Value = 1.25

# --- Row 17: new year boundary marker "2024" (text, matches A10's "2023" look) ---
#     (set first so the shared-string table gets "2024" before the row 16 strings)
$ws.Range("A10").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "'2024"

# --- Row 16: Dec 2023 sub-total row (VL 6-0-0 leave taken) ---
$ws.Range("A16").Value = 45291
$ws.Range("B16").Value = "VL(6-0-0)"
$ws.Range("D16").Value = 6
$ws.Range("K16").Value = "12/12-17/2023"

# --- Rows 18-20: dates for Jan-Mar 2024 ---
$ws.Range("A18").Value = 45322
$ws.Range("A19").Value = 45351
$ws.Range("A20").Value = 45382

# --- Shift the "year boundary" style marker (s=48) down by one row within
#     each subsequent 13-row block, since the new year (2024) now starts
#     one row earlier relative to the block (row 17 instead of the block's
#     usual row). ---
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A75").Copy()
$ws.Range("A76").PasteSpecial(-4122)

# now restore the plain style on the rows that used to carry the marker
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A35").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$ws.Range("A61").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A74").Copy()
$ws.Range("A75").PasteSpecial(-4122)

# --- Extend the leave table by one row (the table's last row carries a
#     heavier "closing" border; push it down and let the newly freed row
#     become a normal entry row). ---
$lo = $ws.ListObjects.Item("Table15")
$ws.Range("A134:K134").Copy($ws.Range("A135:K135"))
$ws.Range("G135").Formula = "=IF(ISBLANK(Table15[[#This Row],[EARNED]]),"""",Table15[[#This Row],[EARNED]])"
$ws.Range("A133:K133").Copy()
$ws.Range("A134:K134").PasteSpecial(-4122)
$lo.Resize($ws.Range("A8:K135"))

# --- CONVERTION sheet: record the VL conversion month (March = 3) ---
$cv.Range("J3").Value = 3
